# Adds a new data row (row 2) to the "avisos" sheet, matching the columns
# Producto | Codigo | Cliente | Telefono | Estado defined in row 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "AURICULAR"
$ws.Range("B2").Value = "b3535"
$ws.Range("C2").Value = "pablo"

# Telefono must stay text (it's a phone number, not a numeric value) so it
# keeps any leading zeros / doesn't get reformatted - force text storage
# before assigning it.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "2494381023"

$ws.Range("E2").Value = "Avisado"
